$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold section headers keep their existing style (s=1); just refresh the text.
$ws.Range("A1").Value = "Heston SL: run_cpu.cpp"
$ws.Range("A2").Value = "Repo: finance.zynqpricer.hls"
$ws.Range("A3").Value = "Version: 052cf806968faa45a7c5b80d785a06e9640806fd"
$ws.Range("A4").Value = "Author: Christian Brugger (brugger@eit.uni-kl.de)"
$ws.Range("A5").Value = "Date: August 30. 2013"

$ws.Range("A7").Value = "Compiler:"
$ws.Range("A8").Value = "Microsoft Visual Studio 2012 Update 3"
$ws.Range("A9").Value = "Flags: /O2 /arch:AVX /fp:fast"
$ws.Range("A10").Value = "Compiled with 32 bit"

$ws.Range("A12").Value = "Code Features:"
$ws.Range("A13").Value = "4 Threads"
$ws.Range("A14").Value = "AVX instructions and vectorized loops"
$ws.Range("A15").Value = "Loop unrolling"
$ws.Range("A16").Value = "Ziguratt transformation (best known method)"

$ws.Range("A18").Value = "Hardware:"
$ws.Range("A19").Value = "Dell Latitude E6430"
$ws.Range("A20").Value = "Intel Core i5-3320M @ 2.60 GHz, 2 cores"

$ws.Range("A22").Value = "Performance:"

# "557568 values / sec" free text becomes a number plus a unit-label cell in B.
$ws.Range("A23").Value = 557568
$ws.Range("B23").Value = "values / sec"

# "= 139.4e6 steps / sec" free text becomes a computed formula plus a unit-label cell.
$ws.Range("A24").Formula = "=A23*250"
$ws.Range("A24").NumberFormat = "0.00E+00"
$ws.Range("B24").Value = "steps / sec"

# Blank spacer row, styled like the scientific-notation cells.
$ws.Range("A25").NumberFormat = "0.00E+00"

# Power section moves down two rows; "35.5 Watt" becomes number + unit-label cell.
$ws.Range("A26").Value = "Power:"
$ws.Range("A27").Value = 35.5
$ws.Range("B27").Value = "Watt"
$ws.Range("A28").Value = "No display"

# Old row 29 (quote-prefixed "= 255 nJ / step") is gone entirely now.
$ws.Range("A29").Clear()

# Computed J/step formula plus its unit-label cell.
$ws.Range("A31").Formula = "=A27/A24"
$ws.Range("A31").NumberFormat = "0.00E+00"
$ws.Range("B31").Value = "J / step"

# New "Power Efficiency:" header, bold like the other section headers.
$ws.Range("A30").Value = "Power Efficiency:"
$ws.Range("A30").Font.Bold = $true

# Trailing blank spacer cell, styled like the old quote-prefixed cells.
$ws.Range("A32").Value = "'x"
$ws.Range("A32").Value = ""

# Column A sized to fit the now mostly-numeric/short content.
$ws.Columns("A:A").AutoFit()

$ws.Range("A29").Select()
